# Applies the pl_mw.xlsx line-results update for the "case with 380 kV" run.
# Updates columns B:F and I:K for rows 2-25 (data rows 0-23) with recalculated
# power-flow line results; columns A, G, H, L:O are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.5701824574025807
$bf[0,1] = 0.2533922926078276
$bf[0,2] = 0.04174904649867273
$bf[0,3] = 0.1189858056445274
$bf[0,4] = 3.500932584304053
$ws.Range("B2:F2").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 1.895790595755926
$ik[0,1] = 0.2297640702139745
$ik[0,2] = 0.8222846735844769
$ws.Range("I2:K2").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.5451197963534753
$bf[0,1] = 0.2438347924552886
$bf[0,2] = 0.04168486166624064
$bf[0,3] = 0.1169706671630273
$bf[0,4] = 3.449881042584849
$ws.Range("B3:F3").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 1.869807077351808
$ik[0,1] = 0.2248288191616368
$ik[0,2] = 0.7878307064502224
$ws.Range("I3:K3").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.5301823612165322
$bf[0,1] = 0.2381351670237564
$bf[0,2] = 0.0416896401398823
$bf[0,3] = 0.1158102007630362
$bf[0,4] = 3.419747242399652
$ws.Range("B4:F4").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 1.854380392277804
$ik[0,1] = 0.2219213151692898
$ik[0,2] = 0.767294596301781
$ws.Range("I4:K4").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.5242086208188823
$bf[0,1] = 0.2358548031491239
$bf[0,2] = 0.04170272962925026
$bf[0,3] = 0.1153566117168836
$bf[0,4] = 3.407771354022088
$ws.Range("B5:F5").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 1.848225932597998
$ik[0,1] = 0.2207672744177174
$ik[0,2] = 0.7590813423154827
$ws.Range("I5:K5").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.5232235325290446
$bf[0,1] = 0.2354786998622558
$bf[0,2] = 0.04170557700557609
$bf[0,3] = 0.1152824596108921
$bf[0,4] = 3.405801094338671
$ws.Range("B6:F6").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 1.847211947982913
$ik[0,1] = 0.2205775042577542
$ik[0,2] = 0.7577269161320999
$ws.Range("I6:K6").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.5301013382417921
$bf[0,1] = 0.2381042422412349
$bf[0,2] = 0.04168977151411113
$bf[0,3] = 0.1158040053285454
$bf[0,4] = 3.419584502192833
$ws.Range("B7:F7").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 1.854296857204346
$ik[0,1] = 0.2219056268028936
$ik[0,2] = 0.7671832006477644
$ws.Range("I7:K7").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.561447171073894
$bf[0,1] = 0.2500617569178587
$bf[0,2] = 0.04171776230855073
$bf[0,3] = 0.1182750306622893
$bf[0,4] = 3.483077886625196
$ws.Range("B8:F8").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 1.88672169365698
$ik[0,1] = 0.2280368672850841
$ik[0,2] = 0.810276374590245
$ws.Range("I8:K8").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.6265030687045225
$bf[0,1] = 0.2748573446426121
$bf[0,2] = 0.0421220540191527
$bf[0,3] = 0.1237313383506446
$bf[0,4] = 3.617255044025484
$ws.Range("B9:F9").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 1.954519826788442
$ik[0,1] = 0.2410387194014305
$ik[0,2] = 0.8997075735611872
$ws.Range("I9:K9").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.6765029606368387
$bf[0,1] = 0.2939095837765535
$bf[0,2] = 0.04263066759814649
$bf[0,3] = 0.128114386141398
$bf[0,4] = 3.721810146115303
$ws.Range("B10:F10").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 2.0069461430285
$ik[0,1] = 0.2511952262824053
$ik[0,2] = 0.9684466452902711
$ws.Range("I10:K10").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.6997317841503445
$bf[0,1] = 0.3027614001306631
$bf[0,2] = 0.04290771334387244
$bf[0,3] = 0.1301901220694219
$bf[0,4] = 3.770690999426989
$ws.Range("B11:F11").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 2.031374494123966
$ik[0,1] = 0.2559485712980774
$ik[0,2] = 1.000384079038383
$ws.Range("I11:K11").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.7085976959401705
$bf[0,1] = 0.3061401482121653
$bf[0,2] = 0.0430191660449637
$bf[0,3] = 0.1309879501351077
$bf[0,4] = 3.789391709831932
$ws.Range("B12:F12").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 2.040708903699013
$ik[0,1] = 0.257767793774434
$ik[0,2] = 1.012574408519555
$ws.Range("I12:K12").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.706685163088423
$bf[0,1] = 0.3054112814927521
$bf[0,2] = 0.04299487219211073
$bf[0,3] = 0.1308155988118216
$bf[0,4] = 3.785355684017986
$ws.Range("B13:F13").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 2.038694829285419
$ik[0,1] = 0.2573751347121771
$ik[0,2] = 1.009944714276003
$ws.Range("I13:K13").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.7004597912516601
$bf[0,1] = 0.3030388344963058
$bf[0,2] = 0.04291675167308284
$bf[0,3] = 0.1302555234199616
$bf[0,4] = 3.772225692414366
$ws.Range("B14:F14").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 2.032140756152955
$ik[0,1] = 0.256097853803638
$ik[0,2] = 1.001385052669804
$ws.Range("I14:K14").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.6966556467000089
$bf[0,1] = 0.3015891320597746
$bf[0,2] = 0.04286975171920204
$bf[0,3] = 0.1299139973311867
$bf[0,4] = 3.76420804503914
$ws.Range("B15:F15").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 2.028137144620445
$ik[0,1] = 0.2553179900153708
$ik[0,2] = 0.9961545652662664
$ws.Range("I15:K15").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.6749946438268353
$bf[0,1] = 0.2933348314472255
$bf[0,2] = 0.0426134788229362
$bf[0,3] = 0.1279803809254716
$bf[0,4] = 3.718642281419534
$ws.Range("B16:F16").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 2.005361403245004
$ik[0,1] = 0.2508872703327114
$ik[0,2] = 0.9663729201319597
$ws.Range("I16:K16").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.6618302671426477
$bf[0,1] = 0.28831855521139
$bf[0,2] = 0.04246794306862967
$bf[0,3] = 0.126815151338505
$bf[0,4] = 3.691027520186765
$ws.Range("B17:F17").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 1.991538045452941
$ik[0,1] = 0.2482033238400021
$ik[0,2] = 0.9482740681223731
$ws.Range("I17:K17").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.6543039791334877
$bf[0,1] = 0.2854507147486629
$bf[0,2] = 0.04238853543444066
$bf[0,3] = 0.1261526464234777
$bf[0,4] = 3.67526825389632
$ws.Range("B18:F18").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 1.983641710682491
$ik[0,1] = 0.246672110981109
$ik[0,2] = 0.93792691389271
$ws.Range("I18:K18").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.6517635208865613
$bf[0,1] = 0.2844826945358534
$bf[0,2] = 0.04236238883715515
$bf[0,3] = 0.1259296562824019
$bf[0,4] = 3.669953705260014
$ws.Range("B19:F19").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 1.98097748836669
$ik[0,1] = 0.2461558164633146
$ik[0,2] = 0.9344343267473505
$ws.Range("I19:K19").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.6632269262556747
$bf[0,1] = 0.2888507454333933
$bf[0,2] = 0.04248299068911621
$bf[0,3] = 0.1269383945088833
$bf[0,4] = 3.693954316266399
$ws.Range("B20:F20").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 1.993003921059824
$ik[0,1] = 0.2484877378486487
$ik[0,2] = 0.9501942173288001
$ws.Range("I20:K20").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.7022864423291821
$bf[0,1] = 0.3037349520892576
$bf[0,2] = 0.04293952023724046
$bf[0,3] = 0.1304197109470948
$bf[0,4] = 3.776077108309607
$ws.Range("B21:F21").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 2.034063563755922
$ik[0,1] = 0.2564724994058309
$ik[0,2] = 1.003896616285687
$ws.Range("I21:K21").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.7282201889411226
$bf[0,1] = 0.313618687605981
$bf[0,2] = 0.0432760056540431
$bf[0,3] = 0.1327636947367594
$bf[0,4] = 3.830860504690577
$ws.Range("B22:F22").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 2.06138786903928
$ik[0,1] = 0.261803157422392
$ik[0,2] = 1.039555794536739
$ws.Range("I22:K22").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.7143416618092431
$bf[0,1] = 0.3083292149435977
$bf[0,2] = 0.04309293793770763
$bf[0,3] = 0.1315063697510155
$bf[0,4] = 3.801519536795752
$ws.Range("B23:F23").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 2.046759384635408
$ik[0,1] = 0.2589477911314049
$ik[0,2] = 1.020472333817793
$ws.Range("I23:K23").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.6625953655829449
$bf[0,1] = 0.2886100921496393
$bf[0,2] = 0.04247617438097251
$bf[0,3] = 0.1268826532204557
$bf[0,4] = 3.692630747922351
$ws.Range("B24:F24").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 1.992341040196621
$ik[0,1] = 0.2483591174550384
$ik[0,2] = 0.949325936871162
$ws.Range("I24:K24").Value = $ik

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 0.6085180375726509
$bf[0,1] = 0.2680038690486413
$bf[0,2] = 0.04197542615521854
$bf[0,3] = 0.1221896828368365
$bf[0,4] = 3.579913227787927
$ws.Range("B25:F25").Value = $bf

$ik = New-Object 'object[,]' 1,3
$ik[0,0] = 1.935722785736957
$ik[0,1] = 0.2374158460400224
$ik[0,2] = 0.8749835650944817
$ws.Range("I25:K25").Value = $ik

